$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '27.296.97'
$ws.Range("E2").Value = '  -4.53%  '
$ws.Range("D3").Value = '1.857.23'
$ws.Range("E3").Value = '  -5.57%  '
$ws.Range("E4").Value = '  -1.28%  '
Set-TextValue "D5" '322.05'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  -1.10%  '
Set-TextValue "D7" '0.4506'
$ws.Range("E7").Value = '  -6.17%  '
Set-TextValue "D8" '0.3850'
$ws.Range("E8").Value = '  -5.40%  '
Set-TextValue "D9" '48.11'
$ws.Range("E9").Value = '  -11.18%  '
Set-TextValue "D10" '0.07872'
$ws.Range("E10").Value = '  -7.62%  '
Set-TextValue "D11" '1.019'
$ws.Range("E11").Value = '  -4.20%  '
Set-TextValue "D12" '21.31'
$ws.Range("E12").Value = '  -5.15%  '
$ws.Range("D13").Value = '1.858.76'
$ws.Range("E13").Value = '  -6.98%  '
Set-TextValue "D14" '7.165'
$ws.Range("E14").Value = '  -6.20%  '
Set-TextValue "D15" '5.869'
$ws.Range("E15").Value = '  -5.38%  '
$ws.Range("E16").Value = '  -1.35%  '
Set-TextValue "D17" '0.00001029'
$ws.Range("E17").Value = '  -4.34%  '
Set-TextValue "D18" '85.35'
$ws.Range("E18").Value = '  -6.48%  '
Set-TextValue "D19" '0.06529'
$ws.Range("E19").Value = '  -1.76%  '
Set-TextValue "D20" '17.00'
$ws.Range("E20").Value = '  -8.79%  '
Set-TextValue "D21" '1.001'
$ws.Range("E21").Value = '  -1.04%  '
Set-TextValue "D22" '5.506'
$ws.Range("E22").Value = '  -6.28%  '
$ws.Range("D23").Value = '27.306.44'
$ws.Range("E23").Value = '  -4.57%  '
Set-TextValue "D24" '10.74'
$ws.Range("E24").Value = '  -6.89%  '
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").Value = '2.084.35'
$ws.Range("E26").Value = '  -6.55%  '
Set-TextValue "D27" '151.72'
$ws.Range("E27").Value = '  -3.08%  '
Set-TextValue "D28" '19.68'
$ws.Range("E28").Value = '  -3.51%  '
Set-TextValue "D29" '2.056'
$ws.Range("E29").Value = '  -5.91%  '
Set-TextValue "D30" '5.468'
$ws.Range("E30").Value = '  -7.49%  '
Set-TextValue "D31" '120.34'
$ws.Range("E31").Value = '  -3.73%  '
Set-TextValue "D32" '1.480'
$ws.Range("E32").Value = '  +1.19%  '
Set-TextValue "D33" '0.09282'
$ws.Range("E33").Value = '  -4.14%  '
Set-TextValue "D34" '0.9317'
$ws.Range("E34").Value = '  -6.11%  '
Set-TextValue "D35" '3.597'
$ws.Range("E35").Value = '  -3.09%  '
Set-TextValue "D36" '5.265'
$ws.Range("E36").Value = '  -6.80%  '
$ws.Range("E37").Value = '  -5.02%  '
Set-TextValue "D38" '0.05973'
$ws.Range("E38").Value = '  -4.45%  '
Set-TextValue "D39" '1.209'
$ws.Range("E39").Value = '  -3.89%  '
Set-TextValue "D40" '8.265'
$ws.Range("E40").Value = '  -9.78%  '
Set-TextValue "D41" '0.9996'
$ws.Range("E41").Value = '  -1.12%  '
Set-TextValue "D42" '0.5905'
$ws.Range("E42").Value = '  -5.48%  '
Set-TextValue "D43" '0.1885'
$ws.Range("E43").Value = '  -2.07%  '
Set-TextValue "D44" '10.12'
Set-TextValue "D45" '1.253'
$ws.Range("E45").Value = '  -8.07%  '
Set-TextValue "D46" '0.5624'
$ws.Range("E46").Value = '  -5.89%  '
$ws.Range("E47").Value = '  -8.59%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D48" '1.920'
$ws.Range("E48").Value = '  -7.47%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D49" '3.352'
$ws.Range("E49").Value = '  -1.93%  '
Set-TextValue "D50" '0.06800'
$ws.Range("E50").Value = '  -0.55%  '
Set-TextValue "D51" '108.02'
$ws.Range("E51").Value = '  -3.25%  '
